# Adicionando valores nos pedidos
# Appends 9 new order rows (41-49) to the pedidos sheet, mirroring the
# columns: id_pedido, id_peca, id_cliente, nome_cliente, id_projeto,
# id_materia_prima, descricao_peca, quantidade, altura_vao, largura_vao,
# altura_peca, largura_peca, area_m2, valor_mp_m2, valor_total, nome_pedido

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: id_pedido, id_peca, id_cliente, nome_cliente, id_projeto, id_materia_prima,
#             descricao_peca, quantidade, altura_vao, largura_vao, altura_peca, largura_peca,
#             area_m2, valor_mp_m2, valor_total, nome_pedido
$data = @(
    @("250213_0001", "250213_0001_001", 1, "DOUGLAS", 8,  11, "Peça Fixa",      2, 1000, 1000, 975,  500, 1,    341.82, 341.82,             "456"),
    @("250213_0001", "250213_0001_002", 1, "DOUGLAS", 8,  11, "Peça Móvel",     2, 1000, 1000, 938,  550, 1.25, 341.82, 427.28,             "456"),
    @("250213_0002", "250213_0002_001", 1, "DOUGLAS", 39, 20, "Peça Principal", 2, 1880, 550,  1880, 550, 2.25, 332.75, 748.6900000000001, "456"),
    @("250213_0003", "250213_0003_001", 1, "DOUGLAS", 26, 11, "Peça Principal", 1, 1845, 600,  1845, 600, 1.25, 341.82, 427.28,             "78945"),
    @("250213_0004", "250213_0004_001", 1, "DOUGLAS", 29, 2,  "Peça Principal", 5, 1845, 750,  1845, 750, 7,    205.75, 1440.25,            "78945"),
    @("250213_0005", "250213_0005_001", 1, "DOUGLAS", 42, 11, "Peça Principal", 3, 1880, 700,  1880, 700, 4,    341.82, 1367.28,            "kijk"),
    @("250213_0006", "250213_0006_001", 1, "DOUGLAS", 51, 2,  "Peça Principal", 1, 938,  450,  938,  450, 0.5,  205.75, 102.88,             "kijk"),
    @("250213_0007", "250213_0007_001", 1, "DOUGLAS", 41, 2,  "Peça Principal", 1, 1880, 650,  1880, 650, 1.25, 205.75, 257.19,             "Boxes casa"),
    @("250213_0008", "250213_0008_001", 1, "DOUGLAS", 42, 2,  "Peça Principal", 1, 1880, 700,  1880, 700, 1.5,  205.75, 308.62,             "Boxes casa")
)

$startRow = 41
# Columns that must stay textual even though their content can look numeric
# (id_pedido/id_peca always, plus nome_pedido which holds values like "456").
$textColumns = @(1, 2, 16)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $rowData = $data[$i]
    for ($col = 1; $col -le $rowData.Length; $col++) {
        $value = $rowData[$col - 1]
        $cell = $ws.Cells.Item($row, $col)
        if ($textColumns -contains $col) {
            # Force text storage (even for purely numeric-looking values like
            # "456" or "78945") without leaving a lasting number-format override
            # on the cell: prefix with an apostrophe so it is entered as text,
            # then reset the cell style back to Normal/default.
            $cell.Value = "'" + $value
            $cell.Style = "Normal"
        } else {
            $cell.Value = $value
        }
    }
}
